$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header columns
$ws.Range("M1").Value = "Color"
$ws.Range("N1").Value = "Talla"

# Add data values for the new columns
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0

# Update the view: scroll so column D is the top-left and select N3
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("N3").Select()
